$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.204320311546326
$ws.Range("B1").Value = 1.384949445724487
$ws.Range("C1").Value = 1.995141983032227
$ws.Range("D1").Value = 1.967057466506958
$ws.Range("E1").Value = 1.012821435928345
